$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated benchmarking results to Hexaly 14: new column headers (B1:I1) ...
$headers = @("acc","T1_U1","T2_U2","T4_U2","T3_U3","T3_U2","T2_U3","T4_U3","T5_U4")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ... and refreshed data values (A2:I21)
$data = @(
    ,(120,0,0,0,0,0,0,0,0)
    ,(60,1,0,0,1,1,1,0,0)
    ,(40,1,1,0,2,2,1,0,0)
    ,(30,1,1,1,2,2,1,1,1)
    ,(24,1,1,2,2,1,1,1,1)
    ,(20,2,2,1,2,3,1,2,2)
    ,(17,2,1,2,2,4,2,2,2)
    ,(15,2,2,1,3,4,2,3,2)
    ,(13,2,4,2,4,3,1,3,2)
    ,(12,2,1,3,3,5,3,3,2)
    ,(10,3,3,3,5,5,3,4,3)
    ,(9,3,3,1,3,9,3,6,3)
    ,(8,3,1,6,5,7,5,4,4)
    ,(7,4,3,6,7,8,5,5,4)
    ,(6,5,6,4,7,10,4,8,5)
    ,(5,5,4,11,10,9,7,6,5)
    ,(4,6,3,17,12,10,9,4,6)
    ,(3,8,7,19,16,14,11,6,8)
    ,(2,10,8,26,13,26,13,6,10)
    ,(1,15,14,55,15,51,20,3,16)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowValues[$c]
    }
}
